# B1--and-B2-PowerPoint.pptx edit
#
# 1) The single table on slide 5 switches its table style (the
#    "tableStyleId" GUID referenced from ppt/tableStyles.xml's builtin
#    gallery) from {1E7C6C0F-1A48-4009-B424-7F911F6028E7} to
#    {FF5AE99B-12BB-46E7-9F9F-074CBE68CEBC}.
#
# 2) The presentation's theme colour palette switches from the
#    "Integral / Red Violet" palette to the default "Office" palette.

$p = $ppt.ActivePresentation

# --- 1) Re-style the table on slide 5 -------------------------------------
$slide = $p.Slides.Item(5)
$tableShape = $slide.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{FF5AE99B-12BB-46E7-9F9F-074CBE68CEBC}")

# --- 2) Swap the theme colour scheme back to the default Office palette ---
$theme = $p.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

# NB: PowerPoint's RGB "Long" packs bytes as 0x00BBGGRR, i.e. R + G*256 + B*65536
$colors.Item(1).RGB  = 0x000000    # dk1      000000
$colors.Item(2).RGB  = 0xFFFFFF    # lt1      FFFFFF
$colors.Item(3).RGB  = 0x6A5444    # dk2      44546A
$colors.Item(4).RGB  = 0xE6E6E7    # lt2      E7E6E6
$colors.Item(5).RGB  = 0xD59B5B    # accent1  5B9BD5
$colors.Item(6).RGB  = 0x317DED    # accent2  ED7D31
$colors.Item(7).RGB  = 0xA5A5A5    # accent3  A5A5A5
$colors.Item(8).RGB  = 0x00C0FF    # accent4  FFC000
$colors.Item(9).RGB  = 0xC47244    # accent5  4472C4
$colors.Item(10).RGB = 0x47AD70    # accent6  70AD47
$colors.Item(11).RGB = 0xC16305    # hlink    0563C1
$colors.Item(12).RGB = 0x724F95    # folHlink 954F72
